# Update BOM workbook: add three new "Diverse" parts to the end of the
# Tabel5 list (rows 52-54), extend the table/dimension accordingly, widen
# column A to fit the new (longer) text, and move the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the "Tabel5" (Diverse) table from A31:B51 to A31:B54 so the three
# new rows become part of the table / autofilter range.
$tbl = $ws.ListObjects.Item("Tabel5")
$tbl.Resize($ws.Range("A31:B54"))

# Row 52 - IEC filtered connector
$ws.Range("A52").Value = "IEC filtered connector Male Schurter, 6A, 250 VAC, Panel mounting, 2 poles"
$ws.Range("B52").Value = "x1"
$ws.Range("C52").Value = "RS"

# Row 53 - Profile VTLB-wire
$ws.Range("A53").Value = "Profile VTLB-wire 3G 0,75mm², white, 2 meter"
$ws.Range("B53").Value = "x1"
$ws.Range("C53").Value = "Hubo"

# Row 54 - Encoder knob
$ws.Range("A54").Value = "Encoder knob"
$ws.Range("B54").Value = "x1"

# Widen column A so the longer descriptions fit (was 53.140625 -> 68.5703125).
$ws.Columns.Item(1).ColumnWidth = 67.59

# Update the visible selection to match the saved view state.
$ws.Range("D48").Select() | Out-Null
